# Update the SmartArt ("Icon Vertical Solid List") text on slide 1.
#   "Introduction to Python, "        -> "Introduction to Python "
#   "Variables/Data Types and "       -> "Variables and Data Types"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape that hosts the SmartArt graphic on this slide.
$smartArtShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasSmartArt) {
        $smartArtShape = $sh
    }
}

$sa = $smartArtShape.SmartArt
$nodes = $sa.AllNodes

for ($i = 1; $i -le $nodes.Count; $i++) {
    $node = $nodes.Item($i)
    $tr = $node.TextFrame2.TextRange
    $t = $tr.Text

    if ($t -eq "Introduction to Python, ") {
        $tr.Text = "Introduction to Python "
    } elseif ($t -eq "Variables/Data Types and ") {
        $tr.Text = "Variables and Data Types"
    }
}
